$wb = $excel.ActiveWorkbook

$newVersion = "mines - January 30 (built on February 02 2026 12.49.33 EST)"

$aboutSheet = $wb.Worksheets.Item("About")
$dataSheet = $wb.Worksheets.Item("Boundaries and methane sources")

$aboutSheet.Range("A2").Value = "Version: $newVersion"
$aboutSheet.Range("A6").Value = "Recommended Citation:  ""Global Energy Monitor, Coal mine boundaries and methane sources for Shaqu No.2 Coal Mine, China, M1195, version '$newVersion'. (See the CC license for attribution requirements if sharing or adapting the data set.)"

for ($row = 2; $row -le 10; $row++) {
    $cell = $dataSheet.Cells.Item($row, 19)
    $cell.Value = $newVersion
}
